$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 78748
$ws.Range("B2").Value = 'Melissa Caldeira'
$ws.Range("C2").Value = 'Operações'
$ws.Range("D2").Value = 'Consulta médica'
$ws.Range("E2").Value = 5
$ws.Range("F2").Value = 45080
$ws.Range("G2").Value = 11403.46

$ws.Range("A3").Value = 84457
$ws.Range("B3").Value = 'Miguel da Cunha'
$ws.Range("C3").Value = 'TI'
$ws.Range("D3").Value = 'Consulta médica'
$ws.Range("E3").Value = 7
$ws.Range("F3").Value = 45094
$ws.Range("G3").Value = 8728.76

$ws.Range("A4").Value = 98998
$ws.Range("B4").Value = 'Ana Sophia Costela'
$ws.Range("C4").Value = 'Jurídico'
$ws.Range("D4").Value = 'Viagem de negócios'
$ws.Range("E4").Value = 7
$ws.Range("F4").Value = 45085
$ws.Range("G4").Value = 5329.3

$ws.Range("A5").Value = 8058
$ws.Range("B5").Value = 'Maysa Gonçalves'
$ws.Range("C5").Value = 'P&D'
$ws.Range("D5").Value = 'Viagem de negócios'
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 45087
$ws.Range("G5").Value = 2817.66

$ws.Range("A6").Value = 65970
$ws.Range("B6").Value = 'Laura Correia'
$ws.Range("C6").Value = 'Vendas'
$ws.Range("D6").Value = 'Consulta médica'
$ws.Range("E6").Value = 5
$ws.Range("F6").Value = 45091
$ws.Range("G6").Value = 4388.12

$ws.Range("A7").Value = 4181
$ws.Range("B7").Value = 'Valentina Martins'
$ws.Range("C7").Value = 'Operações'
$ws.Range("D7").Value = 'Problemas pessoais'
$ws.Range("E7").Value = 8
$ws.Range("F7").Value = 45097
$ws.Range("G7").Value = 2678.12

$ws.Range("A8").Value = 84337
$ws.Range("B8").Value = 'Lara da Cunha'
$ws.Range("C8").Value = 'Atendimento ao Cliente'
$ws.Range("D8").Value = 'Problemas pessoais'
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 45081
$ws.Range("G8").Value = 11302

$ws.Range("A9").Value = 22653
$ws.Range("B9").Value = 'Ana Ramos'
$ws.Range("C9").Value = 'Financeiro'
$ws.Range("D9").Value = 'Viagem de negócios'
$ws.Range("E9").Value = 7
$ws.Range("F9").Value = 45083
$ws.Range("G9").Value = 6423.26

$ws.Range("A10").Value = 72909
$ws.Range("B10").Value = 'Brenda Aragão'
$ws.Range("C10").Value = 'Recursos Humanos'
$ws.Range("D10").Value = 'Problemas pessoais'
$ws.Range("E10").Value = 5
$ws.Range("F10").Value = 45105
$ws.Range("G10").Value = 5099.32

$ws.Range("A11").Value = 6888
$ws.Range("B11").Value = 'Lucca Souza'
$ws.Range("C11").Value = 'Atendimento ao Cliente'
$ws.Range("D11").Value = 'Outros'
$ws.Range("E11").Value = 7
$ws.Range("F11").Value = 45104
$ws.Range("G11").Value = 12245.6

